$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 545.25
$ws.Cells.Item(2, 9).Value = 593.6667
$ws.Cells.Item(2, 10).Value = 400
$ws.Cells.Item(2, 11).Value = 593.6667
$ws.Cells.Item(2, 12).Value = 400
$ws.Cells.Item(2, 13).Value = -480.6667
$ws.Cells.Item(2, 14).Value = -626
$ws.Cells.Item(4, 8).Value = 101.57143
$ws.Cells.Item(4, 9).Value = 101.833336
$ws.Cells.Item(4, 11).Value = 101.833336
$ws.Cells.Item(4, 13).Value = 12.166664
$ws.Cells.Item(8, 8).Value = 742.6
$ws.Cells.Item(9, 8).Value = 2200
$ws.Cells.Item(9, 9).Value = 250
$ws.Cells.Item(9, 10).Value = 10000
$ws.Cells.Item(9, 11).Value = 250
$ws.Cells.Item(9, 12).Value = 10000
$ws.Cells.Item(9, 13).Value = -81
$ws.Cells.Item(9, 14).Value = -10338
$ws.Cells.Item(62, 8).Value = 2250.5
$ws.Cells.Item(62, 9).Value = 2300.8
$ws.Cells.Item(62, 10).Value = 1999
$ws.Cells.Item(62, 11).Value = 2300.8
$ws.Cells.Item(62, 12).Value = 1999
$ws.Cells.Item(62, 13).Value = -1676.8
$ws.Cells.Item(62, 14).Value = -3247
$ws.Cells.Item(65, 8).Value = 2250.5
$ws.Cells.Item(65, 9).Value = 2300.8
$ws.Cells.Item(65, 10).Value = 1999
$ws.Cells.Item(65, 11).Value = 11504
$ws.Cells.Item(65, 12).Value = 9995
$ws.Cells.Item(65, 13).Value = -8384
$ws.Cells.Item(65, 14).Value = -16235
$ws.Cells.Item(86, 8).Value = 2466.9333
$ws.Cells.Item(86, 9).Value = 870
$ws.Cells.Item(86, 10).Value = 5660.8
$ws.Cells.Item(86, 11).Value = 870
$ws.Cells.Item(86, 12).Value = 5660.8
$ws.Cells.Item(86, 13).Value = 253
$ws.Cells.Item(86, 14).Value = -7906.8
$ws.Cells.Item(89, 8).Value = 2466.9333
$ws.Cells.Item(89, 9).Value = 870
$ws.Cells.Item(89, 10).Value = 5660.8
$ws.Cells.Item(89, 11).Value = 4350
$ws.Cells.Item(89, 12).Value = 28304
$ws.Cells.Item(89, 13).Value = 1266
$ws.Cells.Item(89, 14).Value = -39536
$ws.Cells.Item(116, 8).Value = 1256587.5
$ws.Cells.Item(116, 9).Value = 2001540
$ws.Cells.Item(116, 10).Value = 15000
$ws.Cells.Item(116, 11).Value = 2001540
$ws.Cells.Item(116, 12).Value = 15000
$ws.Cells.Item(116, 13).Value = -1998098
$ws.Cells.Item(116, 14).Value = -21884
$ws.Cells.Item(132, 8).Value = 35859740
$ws.Cells.Item(132, 9).Value = 38617530
$ws.Cells.Item(132, 11).Value = 115852590
$ws.Cells.Item(132, 13).Value = -115850060
$ws.Cells.Item(138, 8).Value = 2859.14
$ws.Cells.Item(138, 9).Value = 1482.4615
$ws.Cells.Item(138, 10).Value = 3342.838
$ws.Cells.Item(138, 11).Value = 4447.3845
$ws.Cells.Item(138, 12).Value = 10028.514
$ws.Cells.Item(138, 13).Value = 692.6154999999999
$ws.Cells.Item(138, 14).Value = -20308.514
$ws.Cells.Item(141, 8).Value = 4995.5938
$ws.Cells.Item(141, 9).Value = 5083.6665
$ws.Cells.Item(141, 10).Value = 4520
$ws.Cells.Item(141, 11).Value = 15250.9995
$ws.Cells.Item(141, 12).Value = 13560
$ws.Cells.Item(141, 13).Value = -10070.9995
$ws.Cells.Item(141, 14).Value = -23920
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 9444.469999999999
$ws.Cells.Item(32, 9).Value = 5803.6523
$ws.Cells.Item(32, 10).Value = 17548.227
$ws.Cells.Item(32, 11).Value = 5803.6523
$ws.Cells.Item(32, 12).Value = 17548.227
$ws.Cells.Item(32, 13).Value = -5516.6523
$ws.Cells.Item(32, 14).Value = -18122.227
$ws.Cells.Item(33, 8).Value = 8000
$ws.Cells.Item(33, 9).Value = 8000
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 8000
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = -7671
$ws.Cells.Item(33, 14).ClearContents()
$ws.Cells.Item(45, 8).Value = 2045.8928
$ws.Cells.Item(45, 9).Value = 907.0833
$ws.Cells.Item(45, 10).Value = 2900
$ws.Cells.Item(45, 11).Value = 907.0833
$ws.Cells.Item(45, 12).Value = 2900
$ws.Cells.Item(45, 13).Value = -530.0833
$ws.Cells.Item(45, 14).Value = -3654
$ws.Cells.Item(61, 8).Value = 2517.2856
$ws.Cells.Item(61, 9).Value = 2076.2
$ws.Cells.Item(61, 10).Value = 3620
$ws.Cells.Item(61, 11).Value = 2076.2
$ws.Cells.Item(61, 12).Value = 3620
$ws.Cells.Item(61, 13).Value = -1864.2
$ws.Cells.Item(61, 14).Value = -4044
$ws.Cells.Item(97, 8).Value = 1498.9231
$ws.Cells.Item(97, 9).Value = 1176
$ws.Cells.Item(97, 10).Value = 3275
$ws.Cells.Item(97, 11).Value = 1176
$ws.Cells.Item(97, 12).Value = 3275
$ws.Cells.Item(97, 13).Value = -680
$ws.Cells.Item(97, 14).Value = -4267
$ws.Cells.Item(112, 8).Value = 26050
$ws.Cells.Item(112, 10).Value = 26050
$ws.Cells.Item(112, 12).Value = 26050
$ws.Cells.Item(112, 14).Value = -29004
$ws.Cells.Item(136, 8).Value = 2517.2856
$ws.Cells.Item(136, 9).Value = 2076.2
$ws.Cells.Item(136, 10).Value = 3620
$ws.Cells.Item(136, 11).Value = 6228.599999999999
$ws.Cells.Item(136, 12).Value = 10860
$ws.Cells.Item(136, 13).Value = -3678.599999999999
$ws.Cells.Item(136, 14).Value = -15960
$ws.Cells.Item(139, 8).Value = 43335.332
$ws.Cells.Item(139, 10).Value = 43335.332
$ws.Cells.Item(139, 12).Value = 43335.332
$ws.Cells.Item(139, 14).Value = -53615.332
$ws.Cells.Item(140, 8).Value = 64757.5
$ws.Cells.Item(140, 10).Value = 64757.5
$ws.Cells.Item(140, 12).Value = 64757.5
$ws.Cells.Item(140, 14).Value = -75117.5
$ws.Cells.Item(141, 8).Value = 102502.25
$ws.Cells.Item(141, 10).Value = 102502.25
$ws.Cells.Item(141, 12).Value = 102502.25
$ws.Cells.Item(141, 14).Value = -112862.25
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(55, 8).Value = 30000
$ws.Cells.Item(55, 10).Value = 30000
$ws.Cells.Item(55, 12).Value = 30000
$ws.Cells.Item(55, 14).Value = -30546
$ws.Cells.Item(122, 8).Value = 42115.5
$ws.Cells.Item(122, 10).Value = 42115.5
$ws.Cells.Item(122, 12).Value = 42115.5
$ws.Cells.Item(122, 14).Value = -51915.5
$ws.Cells.Item(134, 8).Value = 4091.5
$ws.Cells.Item(134, 9).Value = 2212.4
$ws.Cells.Item(134, 11).Value = 6637.200000000001
$ws.Cells.Item(134, 13).Value = -4102.200000000001
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(113, 8).Value = 756.725
$ws.Cells.Item(113, 9).Value = 658.1786
$ws.Cells.Item(113, 10).Value = 986.6667
$ws.Cells.Item(113, 11).Value = 1974.5358
$ws.Cells.Item(113, 12).Value = 2960.0001
$ws.Cells.Item(113, 13).Value = 195.4642000000001
$ws.Cells.Item(113, 14).Value = -7300.0001
$ws.Cells.Item(121, 8).Value = 2231.6924
$ws.Cells.Item(121, 9).Value = 332.25
$ws.Cells.Item(121, 10).Value = 2389.9792
$ws.Cells.Item(121, 11).Value = 996.75
$ws.Cells.Item(121, 12).Value = 7169.937600000001
$ws.Cells.Item(121, 13).Value = 313.25
$ws.Cells.Item(121, 14).Value = -9789.937600000001
$ws.Cells.Item(122, 8).Value = 2478.7144
$ws.Cells.Item(122, 9).Value = 534.2381
$ws.Cells.Item(122, 10).Value = 3450.9524
$ws.Cells.Item(122, 11).Value = 4808.142900000001
$ws.Cells.Item(122, 12).Value = 31058.5716
$ws.Cells.Item(122, 13).Value = -2358.142900000001
$ws.Cells.Item(122, 14).Value = -35958.5716
$ws.Cells.Item(132, 8).Value = 4185.3184
$ws.Cells.Item(132, 9).Value = 1151.6923
$ws.Cells.Item(132, 10).Value = 8567.223
$ws.Cells.Item(132, 11).Value = 10365.2307
$ws.Cells.Item(132, 12).Value = 77105.007
$ws.Cells.Item(132, 13).Value = -7835.2307
$ws.Cells.Item(132, 14).Value = -82165.007
$ws.Cells.Item(140, 8).Value = 43894.082
$ws.Cells.Item(140, 9).Value = 73104.14
$ws.Cells.Item(140, 11).Value = 219312.42
$ws.Cells.Item(140, 13).Value = -214132.42
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(93, 8).Value = 18249.875
$ws.Cells.Item(93, 10).Value = 18249.875
$ws.Cells.Item(93, 12).Value = 18249.875
$ws.Cells.Item(93, 14).Value = -21993.875
$ws.Cells.Item(102, 8).Value = 3024.5334
$ws.Cells.Item(102, 9).Value = 2548.7827
$ws.Cells.Item(102, 10).Value = 4587.7144
$ws.Cells.Item(102, 11).Value = 2548.7827
$ws.Cells.Item(102, 12).Value = 4587.7144
$ws.Cells.Item(102, 13).Value = -926.7827000000002
$ws.Cells.Item(102, 14).Value = -7831.7144
$ws.Cells.Item(122, 8).Value = 6881.3335
$ws.Cells.Item(122, 9).Value = 1200
$ws.Cells.Item(122, 10).Value = 9722
$ws.Cells.Item(122, 11).Value = 3600
$ws.Cells.Item(122, 12).Value = 29166
$ws.Cells.Item(122, 13).Value = -1150
$ws.Cells.Item(122, 14).Value = -34066
$ws.Cells.Item(132, 8).Value = 3846.4119
$ws.Cells.Item(132, 9).Value = 1955.5
$ws.Cells.Item(132, 10).Value = 4098.533
$ws.Cells.Item(132, 11).Value = 5866.5
$ws.Cells.Item(132, 12).Value = 12295.599
$ws.Cells.Item(132, 13).Value = -3336.5
$ws.Cells.Item(132, 14).Value = -17355.599
$ws.Cells.Item(138, 8).Value = 52196.668
$ws.Cells.Item(138, 10).Value = 52196.668
$ws.Cells.Item(138, 12).Value = 52196.668
$ws.Cells.Item(138, 14).Value = -62476.668
$ws.Cells.Item(139, 8).Value = 69993.336
$ws.Cells.Item(139, 10).Value = 69993.336
$ws.Cells.Item(139, 12).Value = 69993.336
$ws.Cells.Item(139, 14).Value = -80273.336
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 5913.7188
$ws.Cells.Item(40, 9).Value = 5433.0386
$ws.Cells.Item(40, 10).Value = 7996.6665
$ws.Cells.Item(40, 11).Value = 5433.0386
$ws.Cells.Item(40, 12).Value = 7996.6665
$ws.Cells.Item(40, 13).Value = -5297.0386
$ws.Cells.Item(40, 14).Value = -8268.666499999999
$ws.Cells.Item(50, 8).Value = 34905.832
$ws.Cells.Item(50, 10).Value = 34905.832
$ws.Cells.Item(50, 12).Value = 34905.832
$ws.Cells.Item(50, 14).Value = -36179.832
$ws.Cells.Item(74, 8).Value = 36526.363
$ws.Cells.Item(74, 9).Value = 11500.5
$ws.Cells.Item(74, 10).Value = 42087.668
$ws.Cells.Item(74, 11).Value = 11500.5
$ws.Cells.Item(74, 12).Value = 42087.668
$ws.Cells.Item(74, 13).Value = -10502.5
$ws.Cells.Item(74, 14).Value = -44083.668
$ws.Cells.Item(77, 8).Value = 36526.363
$ws.Cells.Item(77, 9).Value = 11500.5
$ws.Cells.Item(77, 10).Value = 42087.668
$ws.Cells.Item(77, 11).Value = 34501.5
$ws.Cells.Item(77, 12).Value = 126263.004
$ws.Cells.Item(77, 13).Value = -29509.5
$ws.Cells.Item(77, 14).Value = -136247.004
$ws.Cells.Item(122, 8).Value = 6587.2856
$ws.Cells.Item(122, 9).Value = 4960.2856
$ws.Cells.Item(122, 10).Value = 8214.286
$ws.Cells.Item(122, 11).Value = 14880.8568
$ws.Cells.Item(122, 12).Value = 24642.858
$ws.Cells.Item(122, 13).Value = -12430.8568
$ws.Cells.Item(122, 14).Value = -29542.858
$ws.Cells.Item(133, 8).Value = 55877.11
$ws.Cells.Item(133, 10).Value = 55877.11
$ws.Cells.Item(133, 12).Value = 55877.11
$ws.Cells.Item(133, 14).Value = -60937.11
$ws.Cells.Item(136, 8).Value = 3547.4333
$ws.Cells.Item(136, 9).Value = 1801.3529
$ws.Cells.Item(136, 10).Value = 5830.769
$ws.Cells.Item(136, 11).Value = 5404.0587
$ws.Cells.Item(136, 12).Value = 17492.307
$ws.Cells.Item(136, 13).Value = -2854.0587
$ws.Cells.Item(136, 14).Value = -22592.307
$ws.Cells.Item(137, 8).Value = 54453
$ws.Cells.Item(137, 9).Value = 40390
$ws.Cells.Item(137, 10).Value = 56015.555
$ws.Cells.Item(137, 11).Value = 40390
$ws.Cells.Item(137, 12).Value = 56015.555
$ws.Cells.Item(137, 13).Value = -35290
$ws.Cells.Item(137, 14).Value = -66215.55499999999
$ws.Cells.Item(138, 8).Value = 159429
$ws.Cells.Item(138, 10).Value = 159429
$ws.Cells.Item(138, 12).Value = 159429
$ws.Cells.Item(138, 14).Value = -169709
$ws.Cells.Item(139, 8).Value = 47778.332
$ws.Cells.Item(139, 10).Value = 49334
$ws.Cells.Item(139, 12).Value = 49334
$ws.Cells.Item(139, 14).Value = -59614
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(75, 8).Value = 27376.666
$ws.Cells.Item(75, 10).Value = 36065
$ws.Cells.Item(75, 12).Value = 36065
$ws.Cells.Item(75, 14).Value = -37937
$ws.Cells.Item(78, 8).Value = 27376.666
$ws.Cells.Item(78, 10).Value = 36065
$ws.Cells.Item(78, 12).Value = 108195
$ws.Cells.Item(78, 14).Value = -117555
